$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234 (shifts existing rows 234..265 down to 235..266)
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new weekly price record
$ws.Range("A234").Value = 10
$ws.Range("B234").Value = "Vega Modelo de Temuco"
$ws.Range("C234").Value = "La Araucanía"
$ws.Range("D234").Value = 45034
$ws.Range("E234").Value = 9
$ws.Range("F234").Value = 100112012
$ws.Range("G234").Value = "Espinaca"
$ws.Range("H234").Value = "Sin especificar"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 35
$ws.Range("K234").Value = 10000
$ws.Range("L234").Value = 10000
$ws.Range("M234").Value = 10000
$ws.Range("N234").Value = "$/docena de atados"
$ws.Range("O234").Value = "Región de La Araucanía"
$ws.Range("P234").Value = 3333
$ws.Range("Q234").Value = 3
$ws.Range("R234").Value = "Hortaliza"
